# Applies the "Trade #8 closed at 2026-02-17 20:48:33" update to the
# live trading results workbook:
#   - Trade #36 (row 37 on "All Trades", row 4 on "MarketMaking") is closed.
#   - A brand-new OPEN trade #69 is appended to both trade logs.
#   - The Summary and Strategy Status roll-up sheets are refreshed to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.28              # Current Capital
$summary.Range("B4").Value = 0.07000000000000001  # Total P&L $
$summary.Range("B5").Value = 0.04                 # Total P&L %
$summary.Range("B6").Value = 36                   # Total Trades
$summary.Range("B8").Value = 16                   # Losing Trades
$summary.Range("B9").Value = 38.89                # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.28   # Capital
$status.Range("D5").Value = 3        # Trades
$status.Range("E5").Value = -0.04    # P&L $
$status.Range("F5").Value = 0.28     # P&L %

# ---------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# 3a. Close out existing trade #36, stored on row 37
$allTrades.Cells.Item(37, 7).Value = 0.08294899999999999   # G - Exit Price
$allTrades.Cells.Item(37, 8).Value = "CLOSED"               # H - Status
$allTrades.Cells.Item(37, 9).Value = -17.0507                # I - P&L %
$allTrades.Cells.Item(37, 10).Value = -0.02                  # J - P&L $
$allTrades.Cells.Item(37, 11).Value = 100.28                 # K - Capital After
$allTrades.Cells.Item(37, 12).Value = "early_exit"           # L - Exit Reason
$allTrades.Cells.Item(37, 13).Value = 0.11                   # M - Duration (min)

# 3b. Append new trade #69 on row 70
$allTrades.Cells.Item(70, 1).Value = 69
$allTrades.Cells.Item(70, 2).Value = "'2026-02-17"
$allTrades.Cells.Item(70, 3).Value = "20:48:26"
$allTrades.Cells.Item(70, 4).Value = "MarketMaking"
$allTrades.Cells.Item(70, 5).Value = "UP"
$allTrades.Cells.Item(70, 6).Value = 0.1
$allTrades.Cells.Item(70, 8).Value = "OPEN"
$allTrades.Cells.Item(70, 9).Value = 0
$allTrades.Cells.Item(70, 10).Value = 0
$allTrades.Cells.Item(70, 11).Value = 100.3009090909091
$allTrades.Cells.Item(70, 13).Value = 0
$allTrades.Cells.Item(70, 14).Value = 0
$allTrades.Cells.Item(70, 15).Value = 0
$allTrades.Cells.Item(70, 16).Value = 0.6
$allTrades.Cells.Item(70, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# 4. MarketMaking sheet
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# 4a. Close out existing trade #36, stored on row 4
$mm.Cells.Item(4, 7).Value = 0.08294899999999999    # G - Exit Price
$mm.Cells.Item(4, 8).Value = "CLOSED"                # H - Status
$mm.Cells.Item(4, 9).Value = -17.0507                 # I - P&L %
$mm.Cells.Item(4, 10).Value = -0.02                   # J - P&L $
$mm.Cells.Item(4, 11).Value = 100.28                  # K - Capital After
$mm.Cells.Item(4, 16).Value = "early_exit"            # P - Exit Reason
$mm.Cells.Item(4, 17).Value = 0.11                    # Q - Duration (min)

# 4b. Append new trade #69 on row 37
$mm.Cells.Item(37, 1).Value = 69
$mm.Cells.Item(37, 2).Value = "'2026-02-17"
$mm.Cells.Item(37, 3).Value = "20:48:26"
$mm.Cells.Item(37, 4).Value = "MarketMaking"
$mm.Cells.Item(37, 5).Value = "UP"
$mm.Cells.Item(37, 6).Value = 0.1
$mm.Cells.Item(37, 8).Value = "OPEN"
$mm.Cells.Item(37, 9).Value = 0
$mm.Cells.Item(37, 10).Value = 0
$mm.Cells.Item(37, 11).Value = 100.3009090909091
$mm.Cells.Item(37, 12).Value = 0
$mm.Cells.Item(37, 13).Value = 0
$mm.Cells.Item(37, 14).Value = 0.6
$mm.Cells.Item(37, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(37, 17).Value = 0
